$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.878.84'
$ws.Range('E2').Value = '  -8.86%  '
$ws.Range('D3').Value = '3.644.04'
$ws.Range('E3').Value = '  -8.62%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.998'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.14%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '559.14'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -8.72%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '169.48'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.40%  '
$ws.Range('D7').Value = '3.634.88'
$ws.Range('E7').Value = '  -8.63%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.614'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -10.66%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.997'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.30%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.691'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -13.42%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.158'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -15.73%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '49.86'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -13.86%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000286'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -15.53%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '10.31'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -12.01%  '
$ws.Range('D15').Value = '4.185.31'
$ws.Range('E15').Value = '  -9.52%  '
$ws.Range('D16').Value = '3.624.78'
$ws.Range('E16').Value = '  -9.22%  '
$ws.Range('E17').Value = '  -3.77%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '18.98'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -9.57%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.58'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -12.86%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.10'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -11.73%  '
$ws.Range('D21').Value = '66.547.50'
$ws.Range('E21').Value = '  -9.22%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '397.70'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -13.57%  '
$ws.Range('E23').Value = '  -11.28%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '86.16'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -10.73%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.97'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -13.45%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '12.42'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -13.60%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.43'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -6.88%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '5.95'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.33%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '3.66'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -13.85%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '9.24'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -13.43%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '31.99'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -12.38%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.49'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -6.35%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '12.26'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -12.62%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '64.19'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -8.60%  '
$ws.Range('E35').Value = '  -12.67%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '42.20'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -14.05%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '579.05'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -9.03%  '
$ws.Range('D38').Value = '0.0₃0878'
$ws.Range('E38').Value = '  -16.53%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.00'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.11%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.997'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.34%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.385'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -11.34%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.130'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -12.39%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.94'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -10.25%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.91'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -14.66%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0426'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -12.74%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.50'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -4.61%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '8.97'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -15.94%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.131'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -12.25%  '
$ws.Range('B49').Value = 'ApeXProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '3.11'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -9.34%  '
$ws.Range('B50').Value = 'WEMIXToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.62'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -5.50%  '
$ws.Range('B51').Value = 'Maker'
$ws.Range('C51').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D51').Value = '2.666.13'
$ws.Range('E51').Value = '  -5.51%  '
